$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'name"'
$ws.Range("A3").Value = "nameD1vNOe[4@"
$ws.Range("A4").Value = "name3Yr 6"
$ws.Range("A5").Value = "namexu"
$ws.Range("A6").Value = "namePD> {x"
